$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching style of existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for new columns I and J, rows 2-11
$iValues = @(1,1,1,1,1,1,1,1,1,6)
$jValues = @(5,5,5,5,3,4,5,5,5,8)

for ($r = 0; $r -lt 10; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
